$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = -18.69197859759818
$ws.Cells.Item(2, 3).Value = 1.857904754497463
$ws.Cells.Item(2, 4).Value = -18.69197859759818
$ws.Cells.Item(2, 5).Value = -18.69197859759818
$ws.Cells.Item(2, 6).Value = -18.69197859759818
$ws.Cells.Item(2, 7).Value = -18.69197859759818
$ws.Cells.Item(2, 8).Value = -18.69197859759818
$ws.Cells.Item(2, 9).Value = -18.69197859759818
$ws.Cells.Item(2, 10).Value = -18.69197859759818
$ws.Cells.Item(2, 11).Value = -18.69197859759818
$ws.Cells.Item(3, 2).Value = -18.69197859759818
$ws.Cells.Item(3, 3).Value = -18.69197859759818
$ws.Cells.Item(3, 4).Value = -18.69197859759818
$ws.Cells.Item(3, 5).Value = -18.69197859759818
$ws.Cells.Item(3, 6).Value = -18.69197859759818
$ws.Cells.Item(3, 7).Value = -18.69197859759818
$ws.Cells.Item(3, 8).Value = -18.69197859759818
$ws.Cells.Item(3, 9).Value = 1.24184026278923
$ws.Cells.Item(3, 10).Value = -18.69197859759818
$ws.Cells.Item(3, 11).Value = -18.69197859759818
$ws.Cells.Item(4, 2).Value = -18.69197859759818
$ws.Cells.Item(4, 3).Value = 2.150768418908052
$ws.Cells.Item(4, 4).Value = -18.69197859759818
$ws.Cells.Item(4, 5).Value = -18.69197859759818
$ws.Cells.Item(4, 6).Value = 3.491339659744876
$ws.Cells.Item(4, 7).Value = -18.69197859759818
$ws.Cells.Item(4, 8).Value = 1.509557873238905
$ws.Cells.Item(4, 9).Value = -18.69197859759818
$ws.Cells.Item(4, 10).Value = 0.9179521900038237
$ws.Cells.Item(4, 11).Value = -18.69197859759818
$ws.Cells.Item(5, 2).Value = -18.69197859759818
$ws.Cells.Item(5, 3).Value = 1.600977678894361
$ws.Cells.Item(5, 4).Value = -18.69197859759818
$ws.Cells.Item(5, 5).Value = -18.69197859759818
$ws.Cells.Item(5, 6).Value = -18.69197859759818
$ws.Cells.Item(5, 7).Value = 2.847671060949335
$ws.Cells.Item(5, 8).Value = -18.69197859759818
$ws.Cells.Item(5, 9).Value = -18.69197859759818
$ws.Cells.Item(5, 10).Value = -18.69197859759818
$ws.Cells.Item(5, 11).Value = -18.69197859759818
$ws.Cells.Item(6, 2).Value = -18.69197859759818
$ws.Cells.Item(6, 3).Value = -18.69197859759818
$ws.Cells.Item(6, 4).Value = -18.69197859759818
$ws.Cells.Item(6, 5).Value = -18.69197859759818
$ws.Cells.Item(6, 6).Value = -18.69197859759818
$ws.Cells.Item(6, 7).Value = -18.69197859759818
$ws.Cells.Item(6, 8).Value = -18.69197859759818
$ws.Cells.Item(6, 9).Value = -18.69197859759818
$ws.Cells.Item(6, 10).Value = -18.69197859759818
$ws.Cells.Item(6, 11).Value = -18.69197859759818
$ws.Cells.Item(7, 2).Value = 2.397384593564865
$ws.Cells.Item(7, 3).Value = -18.69197859759818
$ws.Cells.Item(7, 4).Value = -18.69197859759818
$ws.Cells.Item(7, 5).Value = -18.69197859759818
$ws.Cells.Item(7, 6).Value = -18.69197859759818
$ws.Cells.Item(7, 7).Value = -18.69197859759818
$ws.Cells.Item(7, 8).Value = -18.69197859759818
$ws.Cells.Item(7, 9).Value = -18.69197859759818
$ws.Cells.Item(7, 10).Value = -18.69197859759818
$ws.Cells.Item(7, 11).Value = -18.69197859759818
$ws.Cells.Item(8, 2).Value = -18.69197859759818
$ws.Cells.Item(8, 3).Value = -18.69197859759818
$ws.Cells.Item(8, 4).Value = -18.69197859759818
$ws.Cells.Item(8, 5).Value = 1.851329873334986
$ws.Cells.Item(8, 6).Value = -18.69197859759818
$ws.Cells.Item(8, 7).Value = -18.69197859759818
$ws.Cells.Item(8, 8).Value = -18.69197859759818
$ws.Cells.Item(8, 9).Value = -18.69197859759818
$ws.Cells.Item(8, 10).Value = -18.69197859759818
$ws.Cells.Item(8, 11).Value = -18.69197859759818
$ws.Cells.Item(9, 2).Value = 3.880831014290335
$ws.Cells.Item(9, 3).Value = -18.69197859759818
$ws.Cells.Item(9, 4).Value = -18.69197859759818
$ws.Cells.Item(9, 5).Value = -18.69197859759818
$ws.Cells.Item(9, 6).Value = -18.69197859759818
$ws.Cells.Item(9, 7).Value = -18.69197859759818
$ws.Cells.Item(9, 8).Value = -18.69197859759818
$ws.Cells.Item(9, 9).Value = -18.69197859759818
$ws.Cells.Item(9, 10).Value = -18.69197859759818
$ws.Cells.Item(9, 11).Value = -18.69197859759818
$ws.Cells.Item(10, 2).Value = -18.69197859759818
$ws.Cells.Item(10, 3).Value = -18.69197859759818
$ws.Cells.Item(10, 4).Value = -18.69197859759818
$ws.Cells.Item(10, 5).Value = -18.69197859759818
$ws.Cells.Item(10, 6).Value = -18.69197859759818
$ws.Cells.Item(10, 7).Value = -18.69197859759818
$ws.Cells.Item(10, 8).Value = -18.69197859759818
$ws.Cells.Item(10, 9).Value = 1.725017829392321
$ws.Cells.Item(10, 10).Value = -18.69197859759818
$ws.Cells.Item(10, 11).Value = 2.214354941861831
$ws.Cells.Item(11, 2).Value = -18.69197859759818
$ws.Cells.Item(11, 3).Value = -18.69197859759818
$ws.Cells.Item(11, 4).Value = -18.69197859759818
$ws.Cells.Item(11, 5).Value = 2.892971520318968
$ws.Cells.Item(11, 6).Value = -18.69197859759818
$ws.Cells.Item(11, 7).Value = 2.845170115646414
$ws.Cells.Item(11, 8).Value = -18.69197859759818
$ws.Cells.Item(11, 9).Value = -18.69197859759818
$ws.Cells.Item(11, 10).Value = -18.69197859759818
$ws.Cells.Item(11, 11).Value = 1.952248256599571
$ws.Cells.Item(12, 2).Value = -18.69197859759818
$ws.Cells.Item(12, 3).Value = -18.69197859759818
$ws.Cells.Item(12, 4).Value = -18.69197859759818
$ws.Cells.Item(12, 5).Value = -18.69197859759818
$ws.Cells.Item(12, 6).Value = -18.69197859759818
$ws.Cells.Item(12, 7).Value = -18.69197859759818
$ws.Cells.Item(12, 8).Value = -18.69197859759818
$ws.Cells.Item(12, 9).Value = -18.69197859759818
$ws.Cells.Item(12, 10).Value = -18.69197859759818
$ws.Cells.Item(12, 11).Value = -18.69197859759818
$ws.Cells.Item(13, 2).Value = -18.69197859759818
$ws.Cells.Item(13, 3).Value = -18.69197859759818
$ws.Cells.Item(13, 4).Value = -18.69197859759818
$ws.Cells.Item(13, 5).Value = 2.629421620367855
$ws.Cells.Item(13, 6).Value = -18.69197859759818
$ws.Cells.Item(13, 7).Value = -18.69197859759818
$ws.Cells.Item(13, 8).Value = -18.69197859759818
$ws.Cells.Item(13, 9).Value = -18.69197859759818
$ws.Cells.Item(13, 10).Value = 1.656646912047303
$ws.Cells.Item(13, 11).Value = 1.758710127031381
$ws.Cells.Item(14, 2).Value = -18.69197859759818
$ws.Cells.Item(14, 3).Value = -18.69197859759818
$ws.Cells.Item(14, 4).Value = -18.69197859759818
$ws.Cells.Item(14, 5).Value = -18.69197859759818
$ws.Cells.Item(14, 6).Value = -18.69197859759818
$ws.Cells.Item(14, 7).Value = -18.69197859759818
$ws.Cells.Item(14, 8).Value = -18.69197859759818
$ws.Cells.Item(14, 9).Value = -18.69197859759818
$ws.Cells.Item(14, 10).Value = -18.69197859759818
$ws.Cells.Item(14, 11).Value = 1.949534304604572
$ws.Cells.Item(15, 2).Value = -18.69197859759818
$ws.Cells.Item(15, 3).Value = -18.69197859759818
$ws.Cells.Item(15, 4).Value = -18.69197859759818
$ws.Cells.Item(15, 5).Value = -18.69197859759818
$ws.Cells.Item(15, 6).Value = -18.69197859759818
$ws.Cells.Item(15, 7).Value = -18.69197859759818
$ws.Cells.Item(15, 8).Value = -18.69197859759818
$ws.Cells.Item(15, 9).Value = -18.69197859759818
$ws.Cells.Item(15, 10).Value = -18.69197859759818
$ws.Cells.Item(15, 11).Value = -18.69197859759818
$ws.Cells.Item(16, 2).Value = -18.69197859759818
$ws.Cells.Item(16, 3).Value = -18.69197859759818
$ws.Cells.Item(16, 4).Value = -18.69197859759818
$ws.Cells.Item(16, 5).Value = -18.69197859759818
$ws.Cells.Item(16, 6).Value = -18.69197859759818
$ws.Cells.Item(16, 7).Value = -18.69197859759818
$ws.Cells.Item(16, 8).Value = -18.69197859759818
$ws.Cells.Item(16, 9).Value = -18.69197859759818
$ws.Cells.Item(16, 10).Value = 1.900847950091344
$ws.Cells.Item(16, 11).Value = -18.69197859759818
$ws.Cells.Item(17, 2).Value = -18.69197859759818
$ws.Cells.Item(17, 3).Value = 2.457032267703853
$ws.Cells.Item(17, 4).Value = -18.69197859759818
$ws.Cells.Item(17, 5).Value = -18.69197859759818
$ws.Cells.Item(17, 6).Value = -18.69197859759818
$ws.Cells.Item(17, 7).Value = -18.69197859759818
$ws.Cells.Item(17, 8).Value = 2.091397768979965
$ws.Cells.Item(17, 9).Value = 2.119704182422165
$ws.Cells.Item(17, 10).Value = 2.557589168239258
$ws.Cells.Item(17, 11).Value = -18.69197859759818
$ws.Cells.Item(18, 2).Value = -18.69197859759818
$ws.Cells.Item(18, 3).Value = -18.69197859759818
$ws.Cells.Item(18, 4).Value = -18.69197859759818
$ws.Cells.Item(18, 5).Value = -18.69197859759818
$ws.Cells.Item(18, 6).Value = -18.69197859759818
$ws.Cells.Item(18, 7).Value = -18.69197859759818
$ws.Cells.Item(18, 8).Value = 1.996767723452342
$ws.Cells.Item(18, 9).Value = 2.056624571438119
$ws.Cells.Item(18, 10).Value = 2.415798473286932
$ws.Cells.Item(18, 11).Value = -18.69197859759818
$ws.Cells.Item(19, 2).Value = -18.69197859759818
$ws.Cells.Item(19, 3).Value = -18.69197859759818
$ws.Cells.Item(19, 4).Value = -18.69197859759818
$ws.Cells.Item(19, 5).Value = -18.69197859759818
$ws.Cells.Item(19, 6).Value = -18.69197859759818
$ws.Cells.Item(19, 7).Value = -18.69197859759818
$ws.Cells.Item(19, 8).Value = 1.633534386348213
$ws.Cells.Item(19, 9).Value = 1.811219216871119
$ws.Cells.Item(19, 10).Value = -18.69197859759818
$ws.Cells.Item(19, 11).Value = -18.69197859759818
$ws.Cells.Item(20, 2).Value = -18.69197859759818
$ws.Cells.Item(20, 3).Value = 0.8955663749755612
$ws.Cells.Item(20, 4).Value = 4.321924858560372
$ws.Cells.Item(20, 5).Value = -18.69197859759818
$ws.Cells.Item(20, 6).Value = 3.129935876525174
$ws.Cells.Item(20, 7).Value = -18.69197859759818
$ws.Cells.Item(20, 8).Value = 1.633402625454653
$ws.Cells.Item(20, 9).Value = 1.210121191860539
$ws.Cells.Item(20, 10).Value = -18.69197859759818
$ws.Cells.Item(20, 11).Value = 2.085084840009241
$ws.Cells.Item(21, 2).Value = -18.69197859759818
$ws.Cells.Item(21, 3).Value = 0.6321074343479779
$ws.Cells.Item(21, 4).Value = -18.69197859759818
$ws.Cells.Item(21, 5).Value = 1.472939374276409
$ws.Cells.Item(21, 6).Value = -18.69197859759818
$ws.Cells.Item(21, 7).Value = 2.489460625181545
$ws.Cells.Item(21, 8).Value = 1.430244019614664
$ws.Cells.Item(21, 9).Value = -18.69197859759818
$ws.Cells.Item(21, 10).Value = -18.69197859759818
$ws.Cells.Item(21, 11).Value = -18.69197859759818
